# Add two new gRNA entries (ARID1A / ARID1B) to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write gene names first (column A), then sequences (column B), so that new
# shared-string entries are created in the same order as the source edit.
$ws.Range("A25").Value = "ARID1A"
$ws.Range("A26").Value = "ARID1B"
$ws.Range("B25").Value = "TCAATCGATGATCTCCCCAT"
$ws.Range("B26").Value = "CCGCAGTACGGACAGCAAGC"

# Copy formatting (font/border/alignment + row height) from the last existing
# data row so the new rows look consistent with the rest of the table.
$ws.Range("A24:B24").Copy()
$ws.Range("A25:B26").PasteSpecial(-4122)
$ws.Rows.Item(25).RowHeight = 17
$ws.Rows.Item(26).RowHeight = 17

# Match the saved selection from the edited workbook.
$ws.Range("L13").Select()
